# feat: add emi to ind_constraint_co2eq row
#
# Adds a new row to the "Process_Set" sheet for the
# "ind_constraint_co2eq" process (id 7) and mirrors the workbook-level
# iterative-calculation tweak that accompanied the change upstream.

$wb = $excel.ActiveWorkbook

# Turn on iterative calculation with a small max-change delta (1E-4),
# matching the updated <calcPr .../> settings in the workbook.
$excel.Iteration = $true
$excel.MaxChange = 0.0001

$ws = $wb.Worksheets.Item("Process_Set")

# New row 8: id = 7, process = "ind_constraint_co2eq"
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "ind_constraint_co2eq"

# Left-align the new label using the same Arial 10 font as the rest
# of the sheet's default style.
$ws.Range("B8").Font.Name = "Arial"
$ws.Range("B8").Font.Size = 10
$ws.Range("B8").HorizontalAlignment = -4131  # xlHAlignLeft

# Leave the active selection where the author ended up after editing.
$ws.Range("B15").Select() | Out-Null
